# Apply the crypto price/volume updates from the Tue Jun 20 04:30:03 UTC 2023
# GitHub Actions automated refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format before writing so that
# numeric-looking strings (e.g. "0.9973", "26.881.55") are preserved
# exactly as text instead of being coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.881.55'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').Value = '1.727.03'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('D4').Value = '0.9973'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '241.85'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').Value = '0.4891'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').Value = '0.2596'
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').Value = '0.06217'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = '1.730.37'
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('D12').Value = '0.06901'
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('D13').Value = '0.6087'
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').Value = '4.486'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').Value = '77.26'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '0.9982'
$ws.Range('D17').Value = '26.639.92'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').Value = '0.9975'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '0.000007173'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').Value = '11.43'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = '1.954.14'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('D22').Value = '4.425'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').Value = '8.565'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').Value = '5.108'
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('D26').Value = '15.33'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('D27').Value = '1.781'
$ws.Range('E27').Value = '  +4.58%  '
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').Value = '3.949'
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').Value = '0.08004'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').Value = '3.688'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').Value = '0.04531'
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('D35').Value = '2.595'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').Value = '0.6251'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').Value = '0.9356'
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('D39').Value = '2.060'
$ws.Range('E39').Value = '  +5.75%  '
$ws.Range('D40').Value = '2.454'
$ws.Range('E40').Value = '  +2.54%  '
$ws.Range('D41').Value = '0.9978'
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('E42').Value = '  +1.22%  '
$ws.Range('D43').Value = '5.651'
$ws.Range('E43').Value = '  +5.69%  '
$ws.Range('D44').Value = '99.34'
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('D45').Value = '0.3856'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').Value = '6.909'
$ws.Range('E46').Value = '  +2.82%  '
$ws.Range('D47').Value = '0.1162'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = '0.05392'
$ws.Range('E48').Value = '  +0.55%  '
$ws.Range('D49').Value = '7.949'
$ws.Range('E49').Value = '  +3.23%  '
$ws.Range('D50').Value = '30.17'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '51.66'
$ws.Range('E51').Value = '  +1.55%  '

# Restore the default "Normal" style so no stray style index is left on
# the cells (matches the original workbook formatting).
$ws.Range("D2:E51").Style = "Normal"

